$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49, pushing the old row 49 (empty separator) and below down by one.
$ws.Rows.Item(49).Insert()

# Row 47: day changed from 4 to 5
$ws.Range("C47").Value = 5

# Row 48: day changed from 4 to 5, end time changed (extends shift)
$ws.Range("C48").Value = 5
$ws.Range("E48").Value = 0.77083333333333337

# New row 49: a new data entry (same day, continuation of work)
$ws.Range("A49").Value = 2014
$ws.Range("B49").Value = 3
$ws.Range("C49").Value = 5
$ws.Range("D49").Value = 0.92708333333333337
$ws.Range("E49").Value = 0.97916666666666663

# Copy styles/format from row 48 D:G into row 49 D:G so the new row matches formatting
$ws.Range("D48:G48").Copy()
$ws.Range("D49:G49").PasteSpecial(-4122)

# Extend the shared formulas in F and G down to the new row 49
$ws.Range("F49").Formula = "=(E49-D49)*24*60"
$ws.Range("G49").Formula = "=F49/60"

# Selection moves to F49
$ws.Range("F49").Select()
